$wb = $excel.ActiveWorkbook

# Update the "Date" metadata value on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-22T16:34:58+00:00"

# Fix the casing of the "Do not know" Display value on the Concepts sheet
# so that it matches the Code column (removes the now-unused "Do not Know" string)
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C8").Value = "Do not know"
